# Adds a new "Hamstring Stretch (Left Leg)" exercise-stats block (columns
# AC:AK) to Sheet1, mirroring the existing per-exercise blocks (A:H, J:R,
# T:AA), and updates the sheet/window view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New block header row (row 1) ---------------------------------------
$ws.Range("AC1").Value2 = "Hamstring Stretch (Left Leg)"

# --- New block sub-headers (row 2) --------------------------------------
$ws.Range("AC2").Value2 = "Rep No"
$ws.Range("AD2").Value2 = "Min Score"
$ws.Range("AE2").Value2 = "Start"
$ws.Range("AF2").Value2 = "Mid"
$ws.Range("AG2").Value2 = "End"
$ws.Range("AH2").Value2 = "Diff 1"
$ws.Range("AI2").Value2 = "Diff 2"
$ws.Range("AJ2").Value2 = "New 1"
$ws.Range("AK2").Value2 = "New 2"

# --- New block data rows (rows 3:11) + per-row diff formulas ------------
$ws.Range("AC3").Value2 = 1
$ws.Range("AD3").Value2 = 0.721581066565861
$ws.Range("AE3").Value2 = 51
$ws.Range("AF3").Value2 = 37
$ws.Range("AG3").Value2 = 4
$ws.Range("AH3").Value2 = 0.0128055810928344
$ws.Range("AI3").Value2 = -0.00756645202636718
$ws.Range("AJ3").Formula = '=AH3-$AH$12'
$ws.Range("AK3").Formula = '=AI3-$AI$12'

$ws.Range("AC4").Value2 = 2
$ws.Range("AD4").Value2 = 0.690776622409206
$ws.Range("AE4").Value2 = 37
$ws.Range("AF4").Value2 = 34
$ws.Range("AG4").Value2 = 4
$ws.Range("AH4").Value2 = -0.0960171222686767
$ws.Range("AI4").Value2 = 0.0355556011199951
$ws.Range("AJ4").Formula = '=AH4-$AH$12'
$ws.Range("AK4").Formula = '=AI4-$AI$12'

$ws.Range("AC5").Value2 = 3
$ws.Range("AD5").Value2 = 0.722689214395245
$ws.Range("AE5").Value2 = 5
$ws.Range("AF5").Value2 = 39
$ws.Range("AG5").Value2 = 3
$ws.Range("AH5").Value2 = 0.0155380964279174
$ws.Range("AI5").Value2 = 0.0340949296951293
$ws.Range("AJ5").Formula = '=AH5-$AH$12'
$ws.Range("AK5").Formula = '=AI5-$AI$12'

$ws.Range("AC6").Value2 = 4
$ws.Range("AD6").Value2 = 0.720705942072139
$ws.Range("AE6").Value2 = 20
$ws.Range("AF6").Value2 = 34
$ws.Range("AG6").Value2 = 3
$ws.Range("AH6").Value2 = -0.00564336776733398
$ws.Range("AI6").Value2 = 0.0588741898536682
$ws.Range("AJ6").Formula = '=AH6-$AH$12'
$ws.Range("AK6").Formula = '=AI6-$AI$12'

$ws.Range("AC7").Value2 = 5
$ws.Range("AD7").Value2 = 0.721190183612289
$ws.Range("AE7").Value2 = 23
$ws.Range("AF7").Value2 = 21
$ws.Range("AG7").Value2 = 3
$ws.Range("AH7").Value2 = 0.0278832912445068
$ws.Range("AI7").Value2 = 0.0294324159622192
$ws.Range("AJ7").Formula = '=AH7-$AH$12'
$ws.Range("AK7").Formula = '=AI7-$AI$12'

$ws.Range("AC8").Value2 = 1
$ws.Range("AD8").Value2 = 0.706382324086765
$ws.Range("AE8").Value2 = 29
$ws.Range("AF8").Value2 = 35
$ws.Range("AG8").Value2 = 2
$ws.Range("AH8").Value2 = -0.0484662055969238
$ws.Range("AI8").Value2 = 0.0821112990379333
$ws.Range("AJ8").Formula = '=AH8-$AH$12'
$ws.Range("AK8").Formula = '=AI8-$AI$12'

$ws.Range("AC9").Value2 = 2
$ws.Range("AD9").Value2 = 0.708468257969004
$ws.Range("AE9").Value2 = 9
$ws.Range("AF9").Value2 = 59
$ws.Range("AG9").Value2 = 4
$ws.Range("AH9").Value2 = -0.0518171787261962
$ws.Range("AI9").Value2 = 0.109801232814788
$ws.Range("AJ9").Formula = '=AH9-$AH$12'
$ws.Range("AK9").Formula = '=AI9-$AI$12'

$ws.Range("AC10").Value2 = 3
$ws.Range("AD10").Value2 = 0.690422420108808
$ws.Range("AE10").Value2 = 8
$ws.Range("AF10").Value2 = 52
$ws.Range("AG10").Value2 = 3
$ws.Range("AH10").Value2 = -0.118322014808654
$ws.Range("AI10").Value2 = 0.142461121082305
$ws.Range("AJ10").Formula = '=AH10-$AH$12'
$ws.Range("AK10").Formula = '=AI10-$AI$12'

$ws.Range("AC11").Value2 = 4
$ws.Range("AD11").Value2 = 0.67636366173635
$ws.Range("AE11").Value2 = 42
$ws.Range("AF11").Value2 = 23
$ws.Range("AG11").Value2 = 3
$ws.Range("AH11").Value2 = -0.204952836036682
$ws.Range("AI11").Value2 = 0.19655704498291
$ws.Range("AJ11").Formula = '=AH11-$AH$12'
$ws.Range("AK11").Formula = '=AI11-$AI$12'

# --- New block averages row (row 12) ------------------------------------
$ws.Range("AD12").Formula = '=AVERAGE(AD3:AD11)'
$ws.Range("AE12").Formula = '=AVERAGE(AE3:AE11)'
$ws.Range("AF12").Formula = '=AVERAGE(AF3:AF11)'
$ws.Range("AG12").Formula = '=AVERAGE(AG3:AG11)'
$ws.Range("AH12").Formula = '=AVERAGE(AH3:AH11)'
$ws.Range("AI12").Formula = '=AVERAGE(AI3:AI11)'
$ws.Range("AJ12").Formula = '=AVERAGE(AJ3:AJ11)'
$ws.Range("AK12").Formula = '=AVERAGE(AK3:AK11)'

# --- View state: scrolled right to column P, zoomed to 130%, and the ----
# --- active selection moved to Z17 --------------------------------------
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 16
$aw.ScrollRow = 1
$aw.Zoom = 130
$ws.Range("Z17").Select()
